# Updated lam and kappa
# Refresh the calibrated parameter values (mu/B/sig2/omega/alpha/beta/gamma/etc.)
# and the associated abs_epsi_autocorr arrays across the AR, SETAR, GARCH,
# TARCH and AR_TARCH calibration sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.0002197522423666759
$ws.Range("B3").Value = 0.7800139040307
$ws.Range("B5").Value = "[0.9999999999999999, 0.25053531886850183, 0.2070088104677465, 0.16140671412818225, 0.2007775953289595, 0.39829783612980957, 0.21785879747660017, 0.21196092783940834, 0.18980277733224551, 0.16636793440461314, 0.16875316747709548, 0.14626408700640145, 0.18955603931193868, 0.18060096593763417, 0.18949999250704327, 0.17724202918153925, 0.1491415219155891, 0.13345667715620596, 0.15272040506725137, 0.16353734030687525]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1260290772802762
$ws.Range("B5").Value = 0.7039927808178683
$ws.Range("B6").Value = 0.08344419462497459
$ws.Range("B7").Value = 0.1375691985345943
$ws.Range("B8").Value = 0.6433513942955118
$ws.Range("B10").Value = "[1.0, 0.24788232362328175, 0.25119466177133426, 0.21728816967872008, 0.20795081477873534, 0.25642478322387396, 0.19650141253355619, 0.17967179125752206, 0.170384936071225, 0.16561474942860355, 0.18641364684488526, 0.16465024152169458, 0.17975112428539128, 0.16064307783176368, 0.16213665194014015, 0.15827317112235512, 0.16320757212572812, 0.15835192559254507, 0.16524687782538058, 0.14445600529099026]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.0004633644944563182
$ws.Range("B3").Value = 0.001250586996403465
$ws.Range("B4").Value = 0.06281857483154474
$ws.Range("B5").Value = 0.9275003932742789
$ws.Range("B6").Value = "[1.0, 0.004201750216821353, -0.02901762303504444, -0.07836316415634041, -0.07385046429193737, 0.29007814442142926, 0.0007028987521514844, 0.03967616936611137, 0.0026336726550008887, -0.035143601374413413, -0.03493933296735276, -0.06959460947950259, -0.004443593538558104, 0.003599233721223633, 0.023253506164405257, 0.024029332489905315, -0.01776286285349463, -0.0535704846833113, -0.06259238877039287, -0.01393493439041259]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.002803422635062306
$ws.Range("B3").Value = 0.001240897985507922
$ws.Range("B4").Value = 0.0525998961843661
$ws.Range("B5").Value = 0.9267060174297441
$ws.Range("B6").Value = "[1.0, 0.002914844080627883, -0.02994861589357286, -0.07922980836185635, -0.07473920778290714, 0.28998149446158655, 0.0002743567606176772, 0.03882068604845218, 0.0018359056840343817, -0.03561130926008774, -0.03575661242165181, -0.06964542166661328, -0.005032276842763942, 0.0029058982549194944, 0.022414491100625874, 0.023342948609641336, -0.018074824292693262, -0.05325493494007477, -0.06263776340128657, -0.01389200897756158]"
$ws.Range("B7").Value = 0.02265231069234847

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.002534223911403063
$ws.Range("B3").Value = 0.0010715122280783
$ws.Range("B4").Value = 0.04380106143789499
$ws.Range("B5").Value = 0.9307266093739496
$ws.Range("B6").Value = "[1.0, 0.017312018758034203, -0.009021291172054724, -0.04495655705498112, -0.042274395089549055, 0.22657016758964318, -0.020382597209890704, 0.017901471982041953, -0.006059452959098886, -0.044405341211522086, -0.019098943552220964, -0.052446148954848185, 0.00023778335427884556, 0.0038200643449998626, 0.007737091260777032, 0.016824270161854497, -0.0227732712133506, -0.04249179335343721, -0.034986202803351095, -0.016552058793695923]"
$ws.Range("B7").Value = 0.03151259845515691
$ws.Range("B9").Value = 0.7843683475058039
